# Append new clock in/out entries to column A, rows 6-16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "OUT -> 2017/02/15 15:03",
    "IN -> 2017/02/15 15:03",
    "OUT -> 2017/02/18 14:33",
    "IN -> 2017/02/18 14:35",
    "OUT -> 2017/02/18 14:35",
    "IN -> 2017/02/18 14:35",
    "OUT -> 2017/02/18 14:37",
    "IN -> 2017/02/18 14:37",
    "OUT -> 2017/02/18 14:38",
    "IN -> 2017/02/18 14:41",
    "OUT -> 2017/02/18 14:41"
)

$row = 6
foreach ($val in $values) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}
